$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '27.208.67'
$ws.Cells.Item(2, 5).Value = '  +0.35%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.851.00'
$ws.Cells.Item(3, 5).Value = '  +0.40%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  -0.13%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '313.36'
$ws.Cells.Item(5, 5).Value = '  +0.26%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -0.27%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.4620'
$ws.Cells.Item(7, 5).Value = '  -0.25%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3690'
$ws.Cells.Item(8, 5).Value = '  -0.45%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.07275'
$ws.Cells.Item(9, 5).Value = '  -1.25%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.8867'
$ws.Cells.Item(10, 5).Value = '  +0.58%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'TRON'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07865'
$ws.Cells.Item(11, 5).Value = '  -0.85%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'Solana'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '19.90'
$ws.Cells.Item(12, 5).Value = '  +0.19%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '5.383'
$ws.Cells.Item(13, 5).Value = '  +0.43%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value = '1.772.88'
$ws.Cells.Item(14, 5).Value = '  -3.19%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Chainlink'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '6.496'
$ws.Cells.Item(15, 5).Value = '  -1.37%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'Litecoin'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '91.53'
$ws.Cells.Item(16, 5).Value = '  -0.24%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '1.003'
$ws.Cells.Item(17, 5).Value = '  -0.29%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.000008850'
$ws.Cells.Item(18, 5).Value = '  -0.85%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.19%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '27.238.06'
$ws.Cells.Item(20, 5).Value = '  +0.32%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '14.66'
$ws.Cells.Item(21, 5).Value = '  -0.95%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '5.066'
$ws.Cells.Item(22, 5).Value = '  -1.29%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'Cosmos'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '10.52'
$ws.Cells.Item(23, 5).Value = '  -0.52%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(24, 4).Value = '2.025.88'
$ws.Cells.Item(24, 5).Value = '  -0.88%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.030'
$ws.Cells.Item(25, 5).Value = '  +8.84%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '151.48'
$ws.Cells.Item(26, 5).Value = '  -0.74%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '18.33'
$ws.Cells.Item(27, 5).Value = '  -0.72%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.030'
$ws.Cells.Item(28, 5).Value = '  -1.79%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '115.67'
$ws.Cells.Item(29, 5).Value = '  -1.18%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '5.001'
$ws.Cells.Item(30, 5).Value = '  -2.32%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.08843'
$ws.Cells.Item(31, 5).Value = '  -0.37%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.138'
$ws.Cells.Item(32, 5).Value = '  +6.06%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.7735'
$ws.Cells.Item(33, 5).Value = '  +4.81%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.504'
$ws.Cells.Item(34, 5).Value = '  +1.02%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.148'
$ws.Cells.Item(35, 5).Value = '  +0.75%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.687'
$ws.Cells.Item(36, 5).Value = '  +8.10%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.104'
$ws.Cells.Item(37, 5).Value = '  +2.24%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.01945'
$ws.Cells.Item(38, 5).Value = '  -0.27%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.05204'
$ws.Cells.Item(39, 5).Value = '  -1.08%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.953'
$ws.Cells.Item(40, 5).Value = '  -0.42%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '6.969'
$ws.Cells.Item(41, 5).Value = '  -1.31%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.5002'
$ws.Cells.Item(42, 5).Value = '  -3.17%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.1614'
$ws.Cells.Item(43, 5).Value = '  -1.17%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '8.411'
$ws.Cells.Item(44, 5).Value = '  +2.79%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.4732'
$ws.Cells.Item(45, 5).Value = '  -2.33%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '10.29'
$ws.Cells.Item(46, 5).Value = '  +0.50%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.26%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '102.26'
$ws.Cells.Item(48, 5).Value = '  -0.18%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.631'
$ws.Cells.Item(49, 5).Value = '  +0.42%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.06193'
$ws.Cells.Item(50, 5).Value = '  -0.45%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '65.30'
$ws.Cells.Item(51, 5).Value = '  -0.10%  '
